$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Previously added"
$ws2 = $wb.Worksheets.Item(2)   # "New"

# ---------------------------------------------------------------------------
# Step 1: move the current "New" rows (2-8) down into "Previously added"
#         as rows 489-495 (values + number formats are preserved by Copy).
# ---------------------------------------------------------------------------
$srcOld = $ws2.Range("A2:F8")
$dstOld = $ws1.Range("A489")
$srcOld.Copy($dstOld)

$oldUrls = @(
  "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/laucesas-pag/icdoc.html",
  "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/graveru-pag/mpcdk.html",
  "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kraslavas-pag/cgjgei.html",
  "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/andrupenes-pag/gdkbg.html",
  "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/udrisu-pag/ngipb.html",
  "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kombulu-pag/pbxmm.html",
  "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kombulu-pag/olnkj.html"
)
for ($i = 0; $i -lt 7; $i++) {
    $row = 489 + $i
    $cell = $ws1.Cells.Item($row, 1)
    $null = $ws1.Hyperlinks.Add($cell, $oldUrls[$i])
}
# Adding a hyperlink re-styles column A with Excel's built-in "Hyperlink"
# style; restore the original (existing) style used by the rest of column A.
$refStyle = $ws1.Cells.Item(488, 1)
$refStyle.Copy()
$ws1.Range("A489:A495").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: remove the now-archived rows from "New" and replace them with the
#         freshly scraped listings.
# ---------------------------------------------------------------------------
$ws2.Hyperlinks.Delete()
$ws2.Rows("6:8").Delete()
$ws2.Range("A2:F5").ClearContents()

# link, price, districtText, areaText, cadastreText, date
$newData = @(
  @("https://www.ss.com/msg/lv/real-estate/wood/aizkraukle-and-reg/kokneses-pag/hmcij.html", "35 000 €", "Aizkraukle un raj.", "7.50 ha.", "32600040206", 46064.67847222222),
  @("https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/zvirgzdenes-pag/bcfjoi.html", "61 500 €", "Ludza un raj.", "21.21 ha.", "68980080265", 46065.509722222225),
  @("https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/kaunatas-pag/obpeg.html", "6 000 €", "Rēzekne un raj.", "2.08 ha.", "78620030347", 46064.75625),
  @("https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/kaunatas-pag/jckep.html", "8 500 €", "Rēzekne un raj.", "2.90 ha.", "78620130045", 46064.74097222222)
)

for ($i = 0; $i -lt 4; $i++) {
    $row = 2 + $i
    $rowData = $newData[$i]

    # Column A first (matches original scrape order -> shared-string order)
    $aCell = $ws2.Cells.Item($row, 1)
    $null = $ws2.Hyperlinks.Add($aCell, $rowData[0])

    $ws2.Cells.Item($row, 2).Value = $rowData[1]
    $ws2.Cells.Item($row, 3).Value = $rowData[2]
    $ws2.Cells.Item($row, 4).Value = $rowData[3]

    # Cadastre numbers must stay text (they look numeric otherwise).
    $eCell = $ws2.Cells.Item($row, 5)
    $eCell.Value = "'" + $rowData[4]

    $ws2.Cells.Item($row, 6).Value = $rowData[5]
}

# Restore original column styles clobbered by the hyperlink add / quote-prefix
# text entry above (column A -> Hyperlink style, column E -> quotePrefix style).
$ws1.Cells.Item(488, 1).Copy()
$ws2.Range("A2:A5").PasteSpecial(-4122)

$ws1.Cells.Item(488, 5).Copy()
$ws2.Range("E2:E5").PasteSpecial(-4122)

Write-Output "done"
